# Commit: "Change ASB to MCSB"
# The workbook contains two worksheets that both show a small block of
# columns labelled "ASB Mapping" / "ASB Guidance" / "ASB Policy" (the
# Azure Security Baseline naming). Microsoft renamed "Azure Security
# Baseline" to "Microsoft Cloud Security Benchmark" (MCSB), so those
# three header labels need to become "MCSB Mapping" / "MCSB Guidance" /
# "MCSB Policy" on every sheet that has them.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Contoso Controls")
$ws2 = $wb.Worksheets.Item("Contoso Controls (Complete)")

# Use Find & Replace across the used cells of each sheet so every
# occurrence of the old "ASB ..." labels becomes the new "MCSB ..." text.
$ws1.Cells.Replace("ASB Mapping", "MCSB Mapping")
$ws1.Cells.Replace("ASB Guidance", "MCSB Guidance")
$ws1.Cells.Replace("ASB Policy", "MCSB Policy")

$ws2.Cells.Replace("ASB Mapping", "MCSB Mapping")
$ws2.Cells.Replace("ASB Guidance", "MCSB Guidance")
$ws2.Cells.Replace("ASB Policy", "MCSB Policy")

# Restore sensible cursor/selection positions on each sheet (leaving the
# second sheet active, as it was originally).
$ws1.Activate()
$ws1.Range("I3").Select()

$ws2.Activate()
$ws2.Range("J7").Select()
